$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "585.04") are stored as strings, matching the inline-string cells
# already used throughout the sheet, instead of being auto-converted to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '63.111.95'
$ws.Range('E2').Value = '  +6.06%  '
$ws.Range('D3').Value = '3.112.75'
$ws.Range('E3').Value = '  +3.96%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '585.04'
$ws.Range('E5').Value = '  +3.80%  '
$ws.Range('D6').Value = '144.35'
$ws.Range('E6').Value = '  +4.03%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.105.97'
$ws.Range('E8').Value = '  +4.14%  '
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('E10').Value = '  +12.83%  '
$ws.Range('D11').Value = '5.79'
$ws.Range('E11').Value = '  +9.97%  '
$ws.Range('E12').Value = '  +3.04%  '
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  +7.86%  '
$ws.Range('D14').Value = '35.53'
$ws.Range('E14').Value = '  +5.09%  '
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '3.632.70'
$ws.Range('E16').Value = '  +4.04%  '
$ws.Range('D17').Value = '7.16'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '63.028.19'
$ws.Range('E18').Value = '  +5.93%  '
$ws.Range('D19').Value = '3.115.52'
$ws.Range('E19').Value = '  +4.17%  '
$ws.Range('D20').Value = '466.34'
$ws.Range('E20').Value = '  +7.37%  '
$ws.Range('D21').Value = '14.07'
$ws.Range('E21').Value = '  +3.86%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  +7.12%  '
$ws.Range('D24').Value = '13.27'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').Value = '81.94'
$ws.Range('E25').Value = '  +2.39%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '8.35'
$ws.Range('E27').Value = '  +7.56%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '2.23'
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('E29').Value = '  +5.46%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '6.84'
$ws.Range('E31').Value = '  +9.67%  '
$ws.Range('D32').Value = '26.89'
$ws.Range('E32').Value = '  +4.44%  '
$ws.Range('D33').Value = '0.109'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').Value = '0.0₃0858'
$ws.Range('E34').Value = '  +10.24%  '
$ws.Range('D35').Value = '2.41'
$ws.Range('E35').Value = '  +15.73%  '
$ws.Range('E36').Value = '  +5.70%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = '6.02'
$ws.Range('E37').Value = '  +2.43%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +18.58%  '
$ws.Range('D39').Value = '50.97'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('D40').Value = '430.44'
$ws.Range('E40').Value = '  +7.33%  '
$ws.Range('D41').Value = '8.71'
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('D42').Value = '2.926.17'
$ws.Range('E42').Value = '  +5.99%  '
$ws.Range('D43').Value = '0.0369'
$ws.Range('E43').Value = '  +4.30%  '
$ws.Range('E44').Value = '  +10.96%  '
$ws.Range('E45').Value = '  +5.48%  '
$ws.Range('E46').Value = '  +7.82%  '
$ws.Range('D47').Value = '35.37'
$ws.Range('E47').Value = '  +3.23%  '
$ws.Range('D49').Value = '123.41'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '24.49'
$ws.Range('E51').Value = '  +4.14%  '

# Restore the default (unstyled) look for column D now that the text values
# are committed, so no cell carries an explicit style that was not there before.
$ws.Range('D2:D51').Style = 'Normal'
